# Apply scheduled market-price + profit updates to the Leve profit sheets.
# Generated from the canonical OOXML diff; each statement targets one cell.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 4921640.5  # H64
$ws.Cells.Item(64, 9).Value = 7356161  # I64
$ws.Cells.Item(64, 10).Value = 52600  # J64
$ws.Cells.Item(64, 11).Value = 7356161  # K64
$ws.Cells.Item(64, 12).Value = 52600  # L64
$ws.Cells.Item(64, 13).Value = -7355913  # M64
$ws.Cells.Item(64, 14).Value = -53096  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 4921640.5  # H67
$ws.Cells.Item(67, 9).Value = 7356161  # I67
$ws.Cells.Item(67, 10).Value = 52600  # J67
$ws.Cells.Item(67, 11).Value = 7356161  # K67
$ws.Cells.Item(67, 12).Value = 52600  # L67
$ws.Cells.Item(67, 13).Value = -7355303  # M67
$ws.Cells.Item(67, 14).Value = -54316  # N67
# Row 106
$ws.Cells.Item(106, 8).Value = 4845.931  # H106
$ws.Cells.Item(106, 9).Value = 5257.136  # I106
$ws.Cells.Item(106, 11).Value = 5257.136  # K106
$ws.Cells.Item(106, 13).Value = -4626.136  # M106
# Row 135
$ws.Cells.Item(135, 8).Value = 2396.9333  # H135
$ws.Cells.Item(135, 9).Value = 1684.4412  # I135
$ws.Cells.Item(135, 11).Value = 15159.9708  # K135
$ws.Cells.Item(135, 13).Value = -12624.9708  # M135
# Row 137
$ws.Cells.Item(137, 8).Value = 1280.6571  # H137
$ws.Cells.Item(137, 9).Value = 1117.0322  # I137
$ws.Cells.Item(137, 11).Value = 3351.0966  # K137
$ws.Cells.Item(137, 13).Value = -801.0966000000003  # M137

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 3858.4583  # H32
$ws.Cells.Item(32, 9).Value = 4079.2559  # I32
$ws.Cells.Item(32, 11).Value = 4079.2559  # K32
$ws.Cells.Item(32, 13).Value = -3792.2559  # M32
# Row 45
$ws.Cells.Item(45, 8).Value = 7974.75  # H45
$ws.Cells.Item(45, 9).Value = 7974.75  # I45
$ws.Cells.Item(45, 11).Value = 7974.75  # K45
$ws.Cells.Item(45, 13).Value = -7597.75  # M45
# Row 97
$ws.Cells.Item(97, 8).Value = 1422.56  # H97
$ws.Cells.Item(97, 9).Value = 678.2  # I97
$ws.Cells.Item(97, 10).Value = 4400  # J97
$ws.Cells.Item(97, 11).Value = 678.2  # K97
$ws.Cells.Item(97, 12).Value = 4400  # L97
$ws.Cells.Item(97, 13).Value = -182.2  # M97
$ws.Cells.Item(97, 14).Value = -5392  # N97
# Row 132
$ws.Cells.Item(132, 8).Value = 5003174  # H132
$ws.Cells.Item(132, 9).Value = 5437998.5  # I132
$ws.Cells.Item(132, 10).Value = 2692  # J132
$ws.Cells.Item(132, 11).Value = 16313995.5  # K132
$ws.Cells.Item(132, 12).Value = 8076  # L132
$ws.Cells.Item(132, 13).Value = -16311465.5  # M132
$ws.Cells.Item(132, 14).Value = -13136  # N132

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Cells.Item(11, 8).Value = 1486.5454  # H11
$ws.Cells.Item(11, 9).Value = 800  # I11
$ws.Cells.Item(11, 10).Value = 1639.1111  # J11
$ws.Cells.Item(11, 11).Value = 800  # K11
$ws.Cells.Item(11, 12).Value = 1639.1111  # L11
$ws.Cells.Item(11, 13).Value = -660  # M11
$ws.Cells.Item(11, 14).Value = -1919.1111  # N11
# Row 105
$ws.Cells.Item(105, 8).Value = 656850.9399999999  # H105
$ws.Cells.Item(105, 9).Value = 1205662.1  # I105
$ws.Cells.Item(105, 11).Value = 1205662.1  # K105
$ws.Cells.Item(105, 13).Value = -1203915.1  # M105
# Row 134
$ws.Cells.Item(134, 8).Value = 1802.6666  # H134
$ws.Cells.Item(134, 9).Value = 1740.9584  # I134
$ws.Cells.Item(134, 11).Value = 5222.8752  # K134
$ws.Cells.Item(134, 13).Value = -2687.8752  # M134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Cells.Item(58, 8).Value = 2298.6365  # H58
$ws.Cells.Item(58, 9).Value = 2009.75  # I58
$ws.Cells.Item(58, 11).Value = 2009.75  # K58
$ws.Cells.Item(58, 13).Value = -1806.75  # M58
# Row 86
$ws.Cells.Item(86, 8).Value = 13718.556  # H86
$ws.Cells.Item(86, 10).Value = 5930  # J86
$ws.Cells.Item(86, 12).Value = 5930  # L86
$ws.Cells.Item(86, 14).Value = -8176  # N86
# Row 89
$ws.Cells.Item(89, 8).Value = 13718.556  # H89
$ws.Cells.Item(89, 10).Value = 5930  # J89
$ws.Cells.Item(89, 12).Value = 29650  # L89
$ws.Cells.Item(89, 14).Value = -40882  # N89
# Row 99
$ws.Cells.Item(99, 8).Value = 6966.2666  # H99
$ws.Cells.Item(99, 10).Value = 18665.334  # J99
$ws.Cells.Item(99, 12).Value = 18665.334  # L99
$ws.Cells.Item(99, 14).Value = -21661.334  # N99
# Row 126
$ws.Cells.Item(126, 8).Value = 6966.2666  # H126
$ws.Cells.Item(126, 10).Value = 18665.334  # J126
$ws.Cells.Item(126, 12).Value = 55996.00199999999  # L126
$ws.Cells.Item(126, 14).Value = -60936.00199999999  # N126
# Row 136
$ws.Cells.Item(136, 8).Value = 2298.6365  # H136
$ws.Cells.Item(136, 9).Value = 2009.75  # I136
$ws.Cells.Item(136, 11).Value = 6029.25  # K136
$ws.Cells.Item(136, 13).Value = -3479.25  # M136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 1761.7778  # H68
$ws.Cells.Item(68, 9).Value = 1375.5  # I68
$ws.Cells.Item(68, 10).Value = 2070.8  # J68
$ws.Cells.Item(68, 11).Value = 4126.5  # K68
$ws.Cells.Item(68, 12).Value = 6212.400000000001  # L68
$ws.Cells.Item(68, 13).Value = -3315.5  # M68
$ws.Cells.Item(68, 14).Value = -7834.400000000001  # N68
# Row 71
$ws.Cells.Item(71, 8).Value = 1761.7778  # H71
$ws.Cells.Item(71, 9).Value = 1375.5  # I71
$ws.Cells.Item(71, 10).Value = 2070.8  # J71
$ws.Cells.Item(71, 11).Value = 12379.5  # K71
$ws.Cells.Item(71, 12).Value = 18637.2  # L71
$ws.Cells.Item(71, 13).Value = -8323.5  # M71
$ws.Cells.Item(71, 14).Value = -26749.2  # N71
# Row 105
$ws.Cells.Item(105, 8).Value = 14909.32  # H105
$ws.Cells.Item(105, 9).Value = 10000  # I105
$ws.Cells.Item(105, 10).Value = 15113.875  # J105
$ws.Cells.Item(105, 11).Value = 30000  # K105
$ws.Cells.Item(105, 12).Value = 45341.625  # L105
$ws.Cells.Item(105, 13).Value = -27379  # M105
$ws.Cells.Item(105, 14).Value = -50583.625  # N105
# Row 113
$ws.Cells.Item(113, 8).Value = 2699.5  # H113
$ws.Cells.Item(113, 10).Value = 2699.5  # J113
$ws.Cells.Item(113, 12).Value = 8098.5  # L113
$ws.Cells.Item(113, 14).Value = -12438.5  # N113
# Row 116
$ws.Cells.Item(116, 8).Value = 14583.25  # H116
$ws.Cells.Item(116, 9).Value = 0  # I116
$ws.Cells.Item(116, 10).Value = 14583.25  # J116
$ws.Cells.Item(116, 11).Value = 0  # K116
$ws.Cells.Item(116, 12).Value = 43749.75  # L116
$ws.Cells.Item(116, 13).ClearContents()  # M116 (removed)
$ws.Cells.Item(116, 14).Value = -50633.75  # N116
# Row 117
$ws.Cells.Item(117, 8).Value = 10328.143  # H117
$ws.Cells.Item(117, 9).Value = 300  # I117
$ws.Cells.Item(117, 10).Value = 14339.4  # J117
$ws.Cells.Item(117, 11).Value = 900  # K117
$ws.Cells.Item(117, 12).Value = 43018.2  # L117
$ws.Cells.Item(117, 13).Value = 2542  # M117
$ws.Cells.Item(117, 14).Value = -49902.2  # N117

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Cells.Item(34, 8).Value = 97499.5  # H34
$ws.Cells.Item(34, 10).Value = 97499.5  # J34
$ws.Cells.Item(34, 12).Value = 97499.5  # L34
$ws.Cells.Item(34, 14).Value = -98035.5  # N34
# Row 64
$ws.Cells.Item(64, 8).Value = 60000.832  # H64
$ws.Cells.Item(64, 10).Value = 60000.832  # J64
$ws.Cells.Item(64, 12).Value = 60000.832  # L64
$ws.Cells.Item(64, 14).Value = -60496.832  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 60000.832  # H67
$ws.Cells.Item(67, 10).Value = 60000.832  # J67
$ws.Cells.Item(67, 12).Value = 60000.832  # L67
$ws.Cells.Item(67, 14).Value = -61716.832  # N67
# Row 76
$ws.Cells.Item(76, 8).Value = 97499.5  # H76
$ws.Cells.Item(76, 10).Value = 97499.5  # J76
$ws.Cells.Item(76, 12).Value = 97499.5  # L76
$ws.Cells.Item(76, 14).Value = -98129.5  # N76
# Row 79
$ws.Cells.Item(79, 8).Value = 97499.5  # H79
$ws.Cells.Item(79, 10).Value = 97499.5  # J79
$ws.Cells.Item(79, 12).Value = 97499.5  # L79
$ws.Cells.Item(79, 14).Value = -99683.5  # N79
# Row 113
$ws.Cells.Item(113, 8).Value = 977583.1  # H113
$ws.Cells.Item(113, 9).Value = 2638.7778  # I113
$ws.Cells.Item(113, 11).Value = 2638.7778  # K113
$ws.Cells.Item(113, 13).Value = -468.7777999999998  # M113
# Row 122
$ws.Cells.Item(122, 8).Value = 1749  # H122
$ws.Cells.Item(122, 9).Value = 1768.1111  # I122
$ws.Cells.Item(122, 11).Value = 5304.3333  # K122
$ws.Cells.Item(122, 13).Value = -2854.3333  # M122
# Row 132
$ws.Cells.Item(132, 8).Value = 1897243.9  # H132
$ws.Cells.Item(132, 9).Value = 3062.3076  # I132
$ws.Cells.Item(132, 11).Value = 9186.9228  # K132
$ws.Cells.Item(132, 13).Value = -6656.9228  # M132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 5119492  # H68
$ws.Cells.Item(68, 9).Value = 6946396  # I68
$ws.Cells.Item(68, 11).Value = 6946396  # K68
$ws.Cells.Item(68, 13).Value = -6945647  # M68
# Row 71
$ws.Cells.Item(71, 8).Value = 5119492  # H71
$ws.Cells.Item(71, 9).Value = 6946396  # I71
$ws.Cells.Item(71, 11).Value = 34731980  # K71
$ws.Cells.Item(71, 13).Value = -34728236  # M71
# Row 122
$ws.Cells.Item(122, 8).Value = 3899.5454  # H122
$ws.Cells.Item(122, 9).Value = 2987  # I122
$ws.Cells.Item(122, 10).Value = 6333  # J122
$ws.Cells.Item(122, 11).Value = 8961  # K122
$ws.Cells.Item(122, 12).Value = 18999  # L122
$ws.Cells.Item(122, 13).Value = -6511  # M122
$ws.Cells.Item(122, 14).Value = -23899  # N122
# Row 123
$ws.Cells.Item(123, 8).Value = 99999  # H123
$ws.Cells.Item(123, 10).Value = 99999  # J123
$ws.Cells.Item(123, 12).Value = 99999  # L123
$ws.Cells.Item(123, 14).Value = -109799  # N123 (new cell)
# Row 132
$ws.Cells.Item(132, 8).Value = 3092.027  # H132
$ws.Cells.Item(132, 9).Value = 2584.963  # I132
$ws.Cells.Item(132, 10).Value = 4461.1  # J132
$ws.Cells.Item(132, 11).Value = 7754.889000000001  # K132
$ws.Cells.Item(132, 12).Value = 13383.3  # L132
$ws.Cells.Item(132, 13).Value = -5224.889000000001  # M132
$ws.Cells.Item(132, 14).Value = -18443.3  # N132
# Row 136
$ws.Cells.Item(136, 8).Value = 10103614  # H136
$ws.Cells.Item(136, 9).Value = 12348096  # I136
$ws.Cells.Item(136, 11).Value = 37044288  # K136
$ws.Cells.Item(136, 13).Value = -37041738  # M136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 6085.1904  # H62
$ws.Cells.Item(62, 9).Value = 3813.4285  # I62
$ws.Cells.Item(62, 11).Value = 3813.4285  # K62
$ws.Cells.Item(62, 13).Value = -3189.4285  # M62
# Row 65
$ws.Cells.Item(65, 8).Value = 6085.1904  # H65
$ws.Cells.Item(65, 9).Value = 3813.4285  # I65
$ws.Cells.Item(65, 11).Value = 19067.1425  # K65
$ws.Cells.Item(65, 13).Value = -15947.1425  # M65
# Row 100
$ws.Cells.Item(100, 8).Value = 553.5294  # H100
$ws.Cells.Item(100, 9).Value = 564.9091  # I100
$ws.Cells.Item(100, 10).Value = 532.6667  # J100
$ws.Cells.Item(100, 11).Value = 1129.8182  # K100
$ws.Cells.Item(100, 12).Value = 1065.3334  # L100
$ws.Cells.Item(100, 13).Value = -588.8181999999999  # M100
$ws.Cells.Item(100, 14).Value = -2147.3334  # N100
# Row 122
$ws.Cells.Item(122, 8).Value = 2313.75  # H122
$ws.Cells.Item(122, 9).Value = 2081.2354  # I122
$ws.Cells.Item(122, 11).Value = 6243.706200000001  # K122
$ws.Cells.Item(122, 13).Value = -3793.706200000001  # M122
# Row 123
$ws.Cells.Item(123, 8).Value = 0  # H123
$ws.Cells.Item(123, 10).Value = 0  # J123
$ws.Cells.Item(123, 12).Value = 0  # L123
$ws.Cells.Item(123, 14).ClearContents()  # N123 (removed)
# Row 124
$ws.Cells.Item(124, 8).Value = 0  # H124
$ws.Cells.Item(124, 10).Value = 0  # J124
$ws.Cells.Item(124, 12).Value = 0  # L124
$ws.Cells.Item(124, 14).ClearContents()  # N124 (removed)
# Row 132
$ws.Cells.Item(132, 8).Value = 2913.258  # H132
$ws.Cells.Item(132, 9).Value = 2674.1155  # I132
$ws.Cells.Item(132, 10).Value = 4156.8  # J132
$ws.Cells.Item(132, 11).Value = 8022.3465  # K132
$ws.Cells.Item(132, 12).Value = 12470.4  # L132
$ws.Cells.Item(132, 13).Value = -5492.3465  # M132
$ws.Cells.Item(132, 14).Value = -17530.4  # N132
